# kicad formating and others
#
# Fill in row 7 (new BOM line: Buzzer) and update the saved view
# (zoom level + active selection) to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New component row: PKM13EPYH4000-A0 / Buzzer ---------------------
$ws.Range("A7").Value = "PKM13EPYH4000-A0"
$ws.Range("B7").Value = "Buzzer"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0.41
$ws.Range("F7").Value = "NO"
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = "N/A"

# --- View state: zoom to 70% and move the active selection ------------
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 70
$ws.Range("D8").Select() | Out-Null
